$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")

# New row 5: "Zibin" / "440807" - added in the same cell style as the rows above
$ws.Range("A5").Value = "Zibin"
$ws.Range("A5").Font.Color = 334100
$ws.Range("B5").Value = "440807"

# Existing numeric totals in column B become text values
$ws.Range("B2").Value = "1"
$ws.Range("B3").Value = "57"
$ws.Range("B4").Value = "420"

# Move the active selection to B4
[void]$ws.Range("B4").Select()
